$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the newly-completed "Real Finalizacion" (actual finish date) values.
#    Copy the date format from an already-formatted cell in the column (I4) so the
#    new cells pick up the existing built-in short-date style instead of creating
#    a brand-new custom number format.
$newDates = @{
    16 = 43810
    19 = 43809
    20 = 43809
    26 = 43810
    27 = 43810
    28 = 43810
    32 = 43810
    33 = 43810
    37 = 43810
}
foreach ($r in $newDates.Keys) {
    $ws.Range("I4").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("I$r").Value = $newDates[$r]
}

# 2. Center (horizontally + vertically) the "Trabajo real (en horas)" and
#    "Duracion Real (en dias)" columns (J:K), matching the rest of the table.
$ws.Columns("J:K").HorizontalAlignment = -4108  # xlCenter
$ws.Columns("J:K").VerticalAlignment = -4108    # xlCenter

# 3. Leave the selection on column K, as in the saved workbook.
$ws.Columns("K:K").Select()
